$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet data lives in an Excel Table ("Tabela1"). Add a new row to
# the table the same way a user/automation (e.g. the daily data-refresh
# bot) would - this keeps the table's `ref`/AutoFilter range, the sheet
# `dimension`, and the sheet data all in sync automatically.
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$newRow = $newListRow.Range

$rowNum = $newRow.Row

# Match the formatting used by the rest of the data rows: column A is a
# date, column B is a thousands-grouped integer, columns C:J are plain
# ("General") numbers. Every data cell also carries the table's thin
# top/left/right rule (no bottom) in the "Calibri Light" 10pt font.
$ws.Cells.Item($rowNum, 1).NumberFormat = "d/\ m/\ yyyy;@"
$ws.Cells.Item($rowNum, 2).NumberFormat = "#,##0"

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item($rowNum, $c)
    $cell.Font.Name = "Calibri Light"
    $cell.Font.Size = 10
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(7).Color = 13998939
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(8).Color = 13998939
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cell.Borders.Item(10).Color = 13998939
    $cell.HorizontalAlignment = -4152     # xlRight
}
$ws.Cells.Item($rowNum, 1).VerticalAlignment = -4160  # xlTop

# New day of Covid-19 data appended to the bottom of the table.
$ws.Cells.Item($rowNum, 1).Value = 43995
$ws.Cells.Item($rowNum, 2).Value = 87386
$ws.Cells.Item($rowNum, 3).Value = 291
$ws.Cells.Item($rowNum, 4).Value = 1495
$ws.Cells.Item($rowNum, 5).Value = 3
$ws.Cells.Item($rowNum, 6).Value = 6
$ws.Cells.Item($rowNum, 7).Value = 0
$ws.Cells.Item($rowNum, 8).Value = 0
$ws.Cells.Item($rowNum, 9).Value = 109
$ws.Cells.Item($rowNum, 10).Value = 0

# Match the saved selection state: the whole new row selected, anchored at A.
$ws.Range("A" + $rowNum + ":J" + $rowNum).Select()
